# "Added plot aggregation to year and parameter set exporting"
#
# Observed workbook-level UI changes:
#   - The lookup sheet "crop sheet lookup" is renamed to "crop_sheet_lookup"
#     (snake_case, matching the other sheet names).
#   - The active/selected tab moves off "farm_layout" and onto the
#     (renamed) "crop_sheet_lookup" sheet, which also gets a new active
#     selection cell.

$wb = $excel.ActiveWorkbook

# Rename "crop sheet lookup" -> "crop_sheet_lookup"
$wsLookup = $wb.Worksheets.Item("crop sheet lookup")
$wsLookup.Name = "crop_sheet_lookup"

# Switch the active/selected sheet from "farm_layout" to the renamed
# "crop_sheet_lookup" sheet, and move its selection to N23.
$wsLookup.Activate()
$wsLookup.Range("N23").Select()
